$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "73.172.99"
$ws.Range("E2").Value = "  +1.51%  "

# Row 3
$ws.Range("D3").Value = "4.050.05"
$ws.Range("E3").Value = "  +0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.10%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "593.68"
$ws.Range("E5").Value = "  +11.66%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.10"
$ws.Range("E6").Value = "  +2.05%  "

# Row 7
$ws.Range("E7").Value = "  -1.66%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("E8").Value = "  -0.04%  "

# Row 9
$ws.Range("E9").Value = "  +1.23%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("E10").Value = "  -0.24%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.85"
$ws.Range("E11").Value = "  +12.76%  "

# Row 12
$ws.Range("E12").Value = "  -0.71%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.03"
$ws.Range("E13").Value = "  +3.75%  "

# Row 14
$ws.Range("D14").Value = "4.704.61"
$ws.Range("E14").Value = "  +0.41%  "

# Row 15
$ws.Range("D15").Value = "4.051.06"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.27"
$ws.Range("E16").Value = "  +1.27%  "

# Row 17
$ws.Range("E17").Value = "  +4.81%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.73"
$ws.Range("E18").Value = "  +0.90%  "

# Row 19
$ws.Range("E19").Value = "  -0.50%  "

# Row 20
$ws.Range("D20").Value = "73.084.49"
$ws.Range("E20").Value = "  +1.41%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "442.46"
$ws.Range("E21").Value = "  +3.03%  "

# Row 22
$ws.Range("E22").Value = "  +12.20%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "97.48"
$ws.Range("E23").Value = "  -0.23%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.52"
$ws.Range("E24").Value = "  +1.78%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.43"
$ws.Range("E25").Value = "  +2.24%  "

# Row 26
$ws.Range("E26").Value = "  +22.24%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.41"
$ws.Range("E27").Value = "  +2.87%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.84"
$ws.Range("E28").Value = "  +1.43%  "

# Row 29
$ws.Range("E29").Value = "  +2.52%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.01"
$ws.Range("E30").Value = "  +0.84%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.01"
$ws.Range("E31").Value = "  +13.53%  "

# Row 32
$ws.Range("E32").Value = "  +4.25%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "13.65"
$ws.Range("E33").Value = "  +2.54%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "691.13"
$ws.Range("E34").Value = "  +2.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "49.10"
$ws.Range("E35").Value = "  +10.41%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "70.88"
$ws.Range("E36").Value = "  +7.43%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.447"
$ws.Range("E37").Value = "  +1.61%  "

# Row 38
$ws.Range("D38").Value = "0.0₃0875"

# Row 39
$ws.Range("B39").Value = "WEMIXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.41"
$ws.Range("E39").Value = "  +7.22%  "

# Row 40
$ws.Range("B40").Value = "THORChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "11.34"
$ws.Range("E40").Value = "  +16.53%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.148"
$ws.Range("E41").Value = "  -2.64%  "

# Row 42
$ws.Range("E42").Value = "  +0.08%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("E43").Value = "  +0.12%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0495"

# Row 45
$ws.Range("E45").Value = "  +0.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.70"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.41"
$ws.Range("E48").Value = "  -0.56%  "

# Row 49
$ws.Range("E49").Value = "  +7.49%  "

# Row 50
$ws.Range("E50").Value = "  +1.72%  "

# Row 51
$ws.Range("E51").Value = "  +8.78%  "
